# Updates symbol list values (price / 1h volume %) and reorders a few
# coin rows (WazirX/One/TigerCash/LEO/GateToken/BTSEToken/BitpandaEcosystemToken)
# to match the refreshed coinranking.com snapshot dated 31-12-2022 20:xx.
#
# Every changed cell in columns D (Price) and E (Volume 1h %) holds text that
# LOOKS numeric ("246.81", "0.71%", ...) but is stored as a plain string in the
# workbook (t="inlineStr"). Writing such a string straight into `.Value` makes
# Excel silently convert it to a real number/percentage and stamp the cell with
# a new number-format style. To keep the cells as plain text with their original
# (unstyled) look, each value is entered with a leading apostrophe (forces text
# entry) and the cell style is then reset to "Normal" to drop the quote-prefix
# style Excel otherwise applies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "246.81"
Set-TextValue "E2" "0.71%"
Set-TextValue "D3" "26.31"
Set-TextValue "E3" "4.80%"
Set-TextValue "D4" "5.074"
Set-TextValue "E4" "1.10%"
Set-TextValue "D7" "0.8132"
Set-TextValue "E7" "0.15%"
Set-TextValue "D8" "0.8440"
Set-TextValue "E8" "0.73%"
Set-TextValue "B9" "One"
Set-TextValue "C9" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D9" "0.0005989"
Set-TextValue "E9" "0.31%"
Set-TextValue "D10" "0.06984"
Set-TextValue "E10" "0.49%"
Set-TextValue "D11" "0.02832"
Set-TextValue "E11" "-0.31%"
Set-TextValue "D12" "0.09416"
Set-TextValue "E12" "0.11%"
Set-TextValue "D13" "0.001524"
Set-TextValue "E13" "0.54%"
Set-TextValue "B14" "TigerCash"
Set-TextValue "C14" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D14" "0.006156"
Set-TextValue "E14" "0.64%"
Set-TextValue "B15" "LEO"
Set-TextValue "C15" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D15" "3.606"
Set-TextValue "E15" "3.10%"
Set-TextValue "B16" "GateToken"
Set-TextValue "C16" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D16" "3.013"
Set-TextValue "E16" "0.17%"
Set-TextValue "B17" "BTSEToken"
Set-TextValue "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D17" "2.056"
Set-TextValue "E17" "-1.71%"
Set-TextValue "B18" "BitpandaEcosystemToken"
Set-TextValue "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D18" "0.3126"
Set-TextValue "E18" "-1.29%"
Set-TextValue "B19" "WazirX"
Set-TextValue "C19" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D19" "0.1338"
Set-TextValue "E19" "-0.05%"
Set-TextValue "D20" "0.03187"
Set-TextValue "E20" "-2.71%"
Set-TextValue "E21" "-1.69%"
Set-TextValue "D22" "3.751"
Set-TextValue "E22" "0.09%"
Set-TextValue "E23" "-0.50%"
Set-TextValue "E24" "-1.48%"
Set-TextValue "D25" "0.001249"
Set-TextValue "E25" "0.64%"
Set-TextValue "D26" "0.004586"
Set-TextValue "D27" "0.00009598"
Set-TextValue "E27" "-1.01%"
Set-TextValue "D28" "0.0001938"
Set-TextValue "E28" "-0.13%"
Set-TextValue "D40" "0.03666"
Set-TextValue "E40" "0.17%"
Set-TextValue "D41" "0.006169"
Set-TextValue "E41" "82.30%"
Set-TextValue "D42" "0.1060"
Set-TextValue "E42" "-21.46%"
Set-TextValue "D44" "0.008257"
Set-TextValue "E44" "2.38%"
Set-TextValue "D45" "0.00005388"
Set-TextValue "E45" "1.79%"
Set-TextValue "E46" "-0.03%"
Set-TextValue "E47" "-38.91%"
Set-TextValue "D48" "0.002598"
Set-TextValue "E48" "27.45%"
Set-TextValue "E49" "-0.03%"
Set-TextValue "E50" "-0.03%"
